$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'27.942.97"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Formula = "'1.634.54"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Formula = "'211.98"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("D6").Formula = "'0.523"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Formula = "'23.40"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -1.38%  '
$ws.Range("E9").Value = '  -1.96%  '
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("D11").Formula = "'0.0881"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +0.59%  '
$ws.Range("D12").Formula = "'1.867.44"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -0.35%  '
$ws.Range("D13").Formula = "'1.637.95"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -0.20%  '
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("E15").Value = '  -1.88%  '
$ws.Range("D16").Formula = "'65.78"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("D17").Formula = "'27.949.14"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").Formula = "'230.86"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Formula = "'0.0₃0725"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Formula = "'7.64"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +0.48%  '
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("D22").Formula = "'10.38"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -7.86%  '
$ws.Range("D23").Formula = "'4.35"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -0.75%  '
$ws.Range("E24").Value = '  +0.26%  '
$ws.Range("D25").Formula = "'155.68"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +2.33%  '
$ws.Range("D26").Formula = "'6.93"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  -0.49%  '
$ws.Range("D28").Formula = "'15.62"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  -0.39%  '
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("D32").Formula = "'3.38"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +1.58%  '
$ws.Range("D33").Formula = "'1.400.57"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -1.70%  '
$ws.Range("E34").Value = '  -1.39%  '
$ws.Range("E37").Value = '  +0.70%  '
$ws.Range("E38").Value = '  +2.02%  '
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("E40").Value = '  -2.68%  '
$ws.Range("E41").Value = '  -0.78%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").Formula = "'66.70"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -0.44%  '
$ws.Range("E44").Value = '  +0.76%  '
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("D47").Formula = "'1.777.29"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("D48").Formula = "'88.16"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("E49").Value = '  -0.74%  '
$ws.Range("D50").Formula = "'0.0998"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("E51").Value = '  -0.31%  '
$ws.Range("B35").Value = 'TrustWalletToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D35").Formula = "'1.05"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +13.81%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Formula = "'1.57"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +0.07%  '
